# test(web)/qa/salesforce/residencial: agregar planes de residencial micronegocio 2p BRM
#
# Update the "Plans" worksheet: refresh the Megas values for plan rows 5-9
# (Residencial / Con_TotalPlay_TV) to the Micronegocio 2P BRM values, and
# clear out row 10 (which previously duplicated the 1000 Megas plan).

$wb = $excel.ActiveWorkbook

$wsPlans = $wb.Worksheets.Item("Plans")

$wsPlans.Range("D5").Value = 50
$wsPlans.Range("D6").Value = 120
$wsPlans.Range("D7").Value = 220
$wsPlans.Range("D8").Value = 520
$wsPlans.Range("D9").Value = 1000

$wsPlans.Range("B10:D10").ClearContents()

# Restore/update the cursor position on each sheet to reflect where the
# edits were made (row 10 of the Plans table).

$wsIntro = $wb.Worksheets.Item("Introduction")
$wsIntro.Range("C8,B10:E10").Select()

$wsTestConfig = $wb.Worksheets.Item("TestConfiguration")
$wsTestConfig.Range("F5,B10:E10").Select()

$wsPlans.Range("B10:E10").Select()

$wsTables = $wb.Worksheets.Item("Tables")
$wsTables.Range("G5,B10:E10").Select()

$wsPlans.Activate()
